$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 1 (title slide)
# ---------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "Understanding Sample Data"
$s1.Shapes.Item(2).TextFrame.TextRange.Text = "A Comprehensive Guide to Analyzing and Interpreting Data"

# ---------------------------------------------------------------
# Slide 2 (Summary / agenda slide)
# ---------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(3).TextFrame.TextRange.Text = "Definition and Types"
$s2.Shapes.Item(4).TextFrame.TextRange.Text = "Importance"
$s2.Shapes.Item(5).TextFrame.TextRange.Text = "Overview of sample data's role in research."
$s2.Shapes.Item(6).TextFrame.TextRange.Text = "Different types of sample data and their significance."
$s2.Shapes.Item(7).TextFrame.TextRange.Text = "Why sample data is essential for research."
$s2.Shapes.Item(8).TextFrame.TextRange.Text = "Recap of key points discussed."
$s2.Shapes.Item(9).TextFrame.TextRange.Text = "Techniques for analyzing sample data."
$s2.Shapes.Item(11).TextFrame.TextRange.Text = "Methods"

# ---------------------------------------------------------------
# Slide 3 (single statement slide)
# ---------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Sample data is crucial for making informed decisions in research and business by providing insights into larger populations."

# ---------------------------------------------------------------
# Slide 4 (Definition and Types)
# ---------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Sample Data"
$s4.Shapes.Item(2).TextFrame.TextRange.Text = "A subset of data collected from a larger population."
$s4.Shapes.Item(3).TextFrame.TextRange.Text = "Includes random, stratified, and systematic samples."
$s4.Shapes.Item(4).TextFrame.TextRange.Text = "Ensures equal chance of selection for all members."
$s4.Shapes.Item(5).TextFrame.TextRange.Text = "Divides population into subgroups for sampling."
$s4.Shapes.Item(6).TextFrame.TextRange.Text = "Definition and Types of Sample Data"
$s4.Shapes.Item(7).TextFrame.TextRange.Text = "Understanding sample data"
$s4.Shapes.Item(8).TextFrame.TextRange.Text = "Types"
$s4.Shapes.Item(9).TextFrame.TextRange.Text = "Random Samples"
$s4.Shapes.Item(10).TextFrame.TextRange.Text = "Stratified Samples"

# ---------------------------------------------------------------
# Slide 5 (Importance of Sample Data in Research)
# ---------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Inferences"
$s5.Shapes.Item(2).TextFrame.TextRange.Text = "Allows researchers to infer about a population."
$s5.Shapes.Item(3).TextFrame.TextRange.Text = "More efficient than studying an entire population."
$s5.Shapes.Item(4).TextFrame.TextRange.Text = "Validates research models effectively."
$s5.Shapes.Item(5).TextFrame.TextRange.Text = "Saves time in data collection and analysis."
$s5.Shapes.Item(6).TextFrame.TextRange.Text = "Importance of Sample Data in Research"
$s5.Shapes.Item(7).TextFrame.TextRange.Text = "Significance in research"
$s5.Shapes.Item(8).TextFrame.TextRange.Text = "Cost-effective"
$s5.Shapes.Item(9).TextFrame.TextRange.Text = "Hypothesis Testing"
$s5.Shapes.Item(10).TextFrame.TextRange.Text = "Time-efficient"

# ---------------------------------------------------------------
# Slide 6 (Methods for Analyzing Sample Data)
# ---------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Descriptive Statistics"
$s6.Shapes.Item(2).TextFrame.TextRange.Text = "Summarizes main features of a data set."
$s6.Shapes.Item(3).TextFrame.TextRange.Text = "Makes predictions about a population."
$s6.Shapes.Item(4).TextFrame.TextRange.Text = "Includes t-tests, chi-square tests, and regression analysis."
$s6.Shapes.Item(5).TextFrame.TextRange.Text = "Aids in decision-making and hypothesis testing."
$s6.Shapes.Item(6).TextFrame.TextRange.Text = "Methods for Analyzing Sample Data"
$s6.Shapes.Item(7).TextFrame.TextRange.Text = "Statistical techniques"
$s6.Shapes.Item(8).TextFrame.TextRange.Text = "Inferential Statistics"
$s6.Shapes.Item(9).TextFrame.TextRange.Text = "Common Methods"
$s6.Shapes.Item(10).TextFrame.TextRange.Text = "Data Analysis"

# ---------------------------------------------------------------
# Slide 7 (Closing slide) - update existing shape text and add the
# remaining placeholders from the "4_TOPIC" layout.
# ---------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Sample Data"

# Switch to the "4_TOPIC" custom layout, which brings in the rest of
# the placeholders (idx 18,19,20,21,11,12,22,23,24) used by the new
# slide content.
$s7.CustomLayout = $p.SlideMaster.CustomLayouts.Item(11)

$s7.Shapes.Item(2).TextFrame.TextRange.Text = "A representative subset used for analysis."
$s7.Shapes.Item(3).TextFrame.TextRange.Text = "Essential for effective research."
$s7.Shapes.Item(4).TextFrame.TextRange.Text = "Various methods employed for analysis."
$s7.Shapes.Item(5).TextFrame.TextRange.Text = "Supports informed decisions in research."
$s7.Shapes.Item(6).TextFrame.TextRange.Text = "Summary"
$s7.Shapes.Item(7).TextFrame.TextRange.Text = "Recap of key ideas"
$s7.Shapes.Item(8).TextFrame.TextRange.Text = "Efficiency"
$s7.Shapes.Item(9).TextFrame.TextRange.Text = "Statistical Methods"
$s7.Shapes.Item(10).TextFrame.TextRange.Text = "Decision-making"

# Match the deck's naming convention ("Text Placeholder <n-1>") for the
# newly added placeholder shapes.
$s7.Shapes.Item(2).Name = "Text Placeholder 2"
$s7.Shapes.Item(3).Name = "Text Placeholder 3"
$s7.Shapes.Item(4).Name = "Text Placeholder 4"
$s7.Shapes.Item(5).Name = "Text Placeholder 5"
$s7.Shapes.Item(6).Name = "Text Placeholder 6"
$s7.Shapes.Item(7).Name = "Text Placeholder 7"
$s7.Shapes.Item(8).Name = "Text Placeholder 8"
$s7.Shapes.Item(9).Name = "Text Placeholder 9"
$s7.Shapes.Item(10).Name = "Text Placeholder 10"
